# Update CDA Logical model for ST.r2b
# Targets the "Metadata" worksheet (sheet1) of the
# StructureDefinition-probability workbook:
#   - bump the Version string
#   - bump the Date timestamp
#   - insert a new "Jurisdiction" property row (empty value) right after
#     "Contact" and before "Description", pushing every row below it down
#     by one

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- simple value updates -------------------------------------------------
$ws.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"

# --- insert the new "Jurisdiction" row ------------------------------------
# Row 11 currently holds "Description"; push it (and everything after it)
# down by one row, then populate the freshly inserted row 11.
$ws.Rows.Item(11).Insert()

# Copy the formatting of the row that is now directly below (the row that
# used to be "Description", i.e. the standard body-row style) onto the new
# row so it matches the rest of the table instead of Excel's blank default.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

$wb.Save()
